# fix: empty cell throw exception
# The "entry=school" row (row 3) used placeholder text instead of being left
# blank like the other "-" cells, which caused a downstream exception on an
# empty cell. Clear C3/D3 back to the normal "-" / base-link placeholder,
# and make sure the "entry=market" link (D4) carries the bold "entry" rich
# text that used to (incorrectly) live on the school row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("school" entry) becomes an empty/placeholder row like the others.
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "os=ios&pf=uc"

# Row 4 ("market" entry) keeps the rich-text formatting (bold "entry") that
# previously lived on the now-blanked school row.
$ws.Range("D4").Value = "os=ios&pf=uc&entry=market"
$boldRun = $ws.Range("D4").Characters(14, 5)
$boldRun.Font.Bold = $true

# Move the active selection to D9, matching the saved view state.
$ws.Range("D9").Select()
